$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.471.99'
$ws.Range("E2").Value = '  -1.18%  '

$ws.Range("D3").Value = '2.986.23'
$ws.Range("E3").Value = '  -2.49%  '

$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").Value = '''501.26'
$ws.Range("E5").Value = '  -2.03%  '

$ws.Range("D6").Value = '''135.64'
$ws.Range("E6").Value = '  +6.42%  '

$ws.Range("E7").Value = '  -0.28%  '

$ws.Range("D8").Value = '''0.427'
$ws.Range("E8").Value = '  -0.76%  '

$ws.Range("D9").Value = '''7.29'
$ws.Range("E9").Value = '  +3.21%  '

$ws.Range("D10").Value = '''0.107'
$ws.Range("E10").Value = '  +3.28%  '

$ws.Range("D11").Value = '''0.354'
$ws.Range("E11").Value = '  -1.49%  '

$ws.Range("E12").Value = '  +0.14%  '

$ws.Range("D13").Value = '3.491.53'
$ws.Range("E13").Value = '  -3.58%  '

$ws.Range("D14").Value = '''25.35'
$ws.Range("E14").Value = '  +4.85%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '''0.0000151'
$ws.Range("E15").Value = '  +5.02%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '56.450.01'
$ws.Range("E16").Value = '  +3.25%  '

$ws.Range("D17").Value = '2.985.88'
$ws.Range("E17").Value = '  -3.36%  '

$ws.Range("D18").Value = '''5.73'
$ws.Range("E18").Value = '  +4.23%  '

$ws.Range("D19").Value = '''12.40'
$ws.Range("E19").Value = '  -0.08%  '

$ws.Range("D20").Value = '''7.83'
$ws.Range("E20").Value = '  +3.60%  '

$ws.Range("D21").Value = '''329.81'
$ws.Range("E21").Value = '  +0.72%  '

$ws.Range("D22").Value = '''1.00'
$ws.Range("E22").Value = '  +0.34%  '

$ws.Range("D23").Value = '''0.473'
$ws.Range("E23").Value = '  -3.07%  '

$ws.Range("D24").Value = '''62.09'
$ws.Range("E24").Value = '  -5.15%  '

$ws.Range("D25").Value = '''0.997'
$ws.Range("E25").Value = '  -0.30%  '

$ws.Range("D26").Value = '''0.164'
$ws.Range("E26").Value = '  -1.26%  '

$ws.Range("D27").Value = '0.0₃0909'
$ws.Range("E27").Value = '  +4.25%  '

$ws.Range("D28").Value = '''0.999'
$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("D29").Value = '''6.48'
$ws.Range("E29").Value = '  -0.27%  '

$ws.Range("D30").Value = '''6.91'
$ws.Range("E30").Value = '  +4.96%  '

$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").Value = '''1.18'
$ws.Range("E31").Value = '  -2.85%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''1.75'
$ws.Range("E32").Value = '  -2.34%  '

$ws.Range("D33").Value = '''20.47'
$ws.Range("E33").Value = '  -0.65%  '

$ws.Range("D34").Value = '''155.35'
$ws.Range("E34").Value = '  -0.50%  '

$ws.Range("D35").Value = '''4.49'
$ws.Range("E35").Value = '  -1.44%  '

$ws.Range("D36").Value = '''1.29'
$ws.Range("E36").Value = '  -0.68%  '

$ws.Range("D37").Value = '''5.61'
$ws.Range("E37").Value = '  -4.45%  '

$ws.Range("D38").Value = '''0.0676'
$ws.Range("E38").Value = '  +3.28%  '

$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D39").Value = '''23.14'
$ws.Range("E39").Value = '  +1.22%  '

$ws.Range("B40").Value = 'RenzoRestakedETH'
$ws.Range("C40").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D40").Value = '3.016.13'
$ws.Range("E40").Value = '  -3.09%  '

$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  -0.11%  '

$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = '''36.31'
$ws.Range("E42").Value = '  +0.32%  '

$ws.Range("D43").Value = '''0.638'
$ws.Range("E43").Value = '  -3.71%  '

$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").Value = '''1.00'
$ws.Range("E44").Value = '  -3.12%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.236.31'
$ws.Range("E45").Value = '  +0.79%  '

$ws.Range("E46").Value = '  +2.65%  '

$ws.Range("D47").Value = '''3.58'
$ws.Range("E47").Value = '  -3.26%  '

$ws.Range("D48").Value = '''1.96'
$ws.Range("E48").Value = '  +14.64%  '

$ws.Range("D49").Value = '''0.0237'
$ws.Range("E49").Value = '  +5.16%  '

$ws.Range("D50").Value = '''5.80'
$ws.Range("E50").Value = '  -2.19%  '

$ws.Range("D51").Value = '''19.17'
$ws.Range("E51").Value = '  -1.71%  '

